# Fruta / hortaliza, semanal
# The weekly data refresh re-sorts/re-orders the existing price rows
# (rows 2-30) of the sheet. Columns A,B,C,E,F,G,H,I,J are constant for
# every row in this sheet, so only the "content" columns D,K,L,M,N,O,P,
# Q,R,S,T need to move to their new row position. The mapping below says,
# for each destination (new) row, which source (old) row's data block it
# receives.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2  = 24
    3  = 17
    4  = 12
    5  = 6
    6  = 20
    7  = 11
    8  = 16
    9  = 23
    10 = 30
    11 = 25
    12 = 19
    13 = 28
    14 = 22
    15 = 10
    16 = 8
    17 = 13
    18 = 26
    19 = 21
    20 = 3
    21 = 4
    22 = 15
    23 = 14
    24 = 18
    25 = 5
    26 = 2
    27 = 29
    28 = 7
    29 = 27
    30 = 9
}

$cols = @("D", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot every source cell's value before writing anything, since the
# mapping is a permutation (sources and destinations overlap).
$snapshot = @{}
foreach ($col in $cols) {
    for ($row = 2; $row -le 30; $row++) {
        $addr = $col + $row
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    foreach ($col in $cols) {
        $srcAddr = $col + $srcRow
        $destAddr = $col + $destRow
        $ws.Range($destAddr).Value2 = $snapshot[$srcAddr]
    }
}
